$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 826.7742
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 826.7742
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2480.3226
$ws.Range("N17").Value = -2816.3226
$ws.Range("H40").Value = 2297.7144
$ws.Range("I40").Value = 1657
$ws.Range("J40").Value = 2472.4546
$ws.Range("K40").Value = 1657
$ws.Range("L40").Value = 2472.4546
$ws.Range("M40").Value = -1482
$ws.Range("N40").Value = -2822.4546
$ws.Range("H103").Value = 3293.375
$ws.Range("I103").Value = 850
$ws.Range("J103").Value = 4759.4
$ws.Range("K103").Value = 2550
$ws.Range("L103").Value = 14278.2
$ws.Range("M103").Value = -1964
$ws.Range("N103").Value = -15450.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1927.8518
$ws.Range("I2").Value = 1812.4762
$ws.Range("J2").Value = 2331.6667
$ws.Range("K2").Value = 1812.4762
$ws.Range("L2").Value = 2331.6667
$ws.Range("M2").Value = -1699.4762
$ws.Range("N2").Value = -2557.6667
$ws.Range("H32").Value = 2866.28
$ws.Range("I32").Value = 2720.0637
$ws.Range("J32").Value = 5157
$ws.Range("K32").Value = 2720.0637
$ws.Range("L32").Value = 5157
$ws.Range("M32").Value = -2433.0637
$ws.Range("N32").Value = -5731
$ws.Range("H88").Value = 2444.182
$ws.Range("I88").Value = 1951.5
$ws.Range("J88").Value = 2725.7144
$ws.Range("K88").Value = 1951.5
$ws.Range("L88").Value = 2725.7144
$ws.Range("M88").Value = -1545.5
$ws.Range("N88").Value = -3537.7144
$ws.Range("H91").Value = 2444.182
$ws.Range("I91").Value = 1951.5
$ws.Range("J91").Value = 2725.7144
$ws.Range("K91").Value = 1951.5
$ws.Range("L91").Value = 2725.7144
$ws.Range("M91").Value = -547.5
$ws.Range("N91").Value = -5533.7144
$ws.Range("H95").Value = 20416.125
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 20416.125
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 20416.125
$ws.Range("N95").Value = -25908.125
$ws.Range("H116").Value = 1927.8518
$ws.Range("I116").Value = 1812.4762
$ws.Range("J116").Value = 2331.6667
$ws.Range("K116").Value = 1812.4762
$ws.Range("L116").Value = 2331.6667
$ws.Range("M116").Value = 481.5237999999999
$ws.Range("N116").Value = -6919.6667
$ws.Range("H122").Value = 1133.1464
$ws.Range("I122").Value = 986.21875
$ws.Range("K122").Value = 2958.65625
$ws.Range("M122").Value = -508.65625
$ws.Range("H132").Value = 5676.116
$ws.Range("I132").Value = 3605.1035
$ws.Range("J132").Value = 9966.071
$ws.Range("K132").Value = 10815.3105
$ws.Range("L132").Value = 29898.213
$ws.Range("M132").Value = -8285.3105
$ws.Range("N132").Value = -34958.213

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1927.8518
$ws.Range("I3").Value = 1812.4762
$ws.Range("J3").Value = 2331.6667
$ws.Range("K3").Value = 1812.4762
$ws.Range("L3").Value = 2331.6667
$ws.Range("M3").Value = -1698.4762
$ws.Range("N3").Value = -2559.6667
$ws.Range("H86").Value = 1587.2222
$ws.Range("I86").Value = 1450
$ws.Range("K86").Value = 1450
$ws.Range("M86").Value = -327
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H89").Value = 1587.2222
$ws.Range("I89").Value = 1450
$ws.Range("K89").Value = 7250
$ws.Range("M89").Value = -1634
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1842.6617
$ws.Range("I58").Value = 1239.0667
$ws.Range("J58").Value = 3023.6086
$ws.Range("K58").Value = 1239.0667
$ws.Range("L58").Value = 3023.6086
$ws.Range("M58").Value = -1036.0667
$ws.Range("N58").Value = -3429.6086
$ws.Range("H88").Value = 28000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 28000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 28000
$ws.Range("N88").Value = -28812
$ws.Range("H91").Value = 28000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 28000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 28000
$ws.Range("N91").Value = -30808
$ws.Range("H136").Value = 1842.6617
$ws.Range("I136").Value = 1239.0667
$ws.Range("J136").Value = 3023.6086
$ws.Range("K136").Value = 3717.2001
$ws.Range("L136").Value = 9070.825800000001
$ws.Range("M136").Value = -1167.2001
$ws.Range("N136").Value = -14170.8258

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5871.4287
$ws.Range("I70").Value = 5100
$ws.Range("J70").Value = 6346.154
$ws.Range("K70").Value = 5100
$ws.Range("L70").Value = 6346.154
$ws.Range("M70").Value = -4830
$ws.Range("N70").Value = -6886.154
$ws.Range("H73").Value = 5871.4287
$ws.Range("I73").Value = 5100
$ws.Range("J73").Value = 6346.154
$ws.Range("K73").Value = 5100
$ws.Range("L73").Value = 6346.154
$ws.Range("M73").Value = -4164
$ws.Range("N73").Value = -8218.154
$ws.Range("H95").Value = 9824.799999999999
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 9824.799999999999
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 9824.799999999999
$ws.Range("N95").Value = -15316.8
$ws.Range("H132").Value = 1070550.9
$ws.Range("I132").Value = 1489599.8
$ws.Range("J132").Value = 3881.0908
$ws.Range("K132").Value = 4468799.4
$ws.Range("L132").Value = 11643.2724
$ws.Range("M132").Value = -4466269.4
$ws.Range("N132").Value = -16703.2724

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 14234.429
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 14234.429
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 14234.429
$ws.Range("N101").Value = -20724.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 5118
$ws.Range("I75").Value = 5118
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 5118
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("M75").Value = -4182
$ws.Range("H78").Value = 5118
$ws.Range("I78").Value = 5118
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 15354
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("M78").Value = -10674
$ws.Range("H98").Value = 37196.668
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 37196.668
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 37196.668
$ws.Range("N98").Value = -43186.668
$ws.Range("H129").Value = 26530
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 26530
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 26530
$ws.Range("N129").Value = -36530
